$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "rank" value for this row (column D, row 5)
$ws.Range("D5").Value = 87

# Add the new "elapsed time" value for this row (column J, row 5)
$ws.Range("J5").Value = 4416

# Update the active selection to D5, matching the saved view state
$ws.Range("D5").Select()
